# Fix the QosProfile gendoc template version/tag string:
#   "...QosProfile_1.0.0-tsp.d.t+gendoc.1..."
# becomes
#   "...QosProfile_1.0.0-tsi.d.t+gendoc.1..."
# (the "p" in "tsp" is replaced by "i", i.e. the pre-release tag is
# corrected to "tsi" instead of "tsp"). This string appears both in the
# body text (inside the <output path='...'/> example) and in the page
# header (document title / version line).

$d = $word.ActiveDocument

# --- main document body -----------------------------------------------
$d.Content.Find.Execute(
    "tsp",   $true, $false, $false, $false, $false,
    $true,   1,     $false, "tsi",  2) | Out-Null

# --- page header ---------------------------------------------------------
$section = $d.Sections.Item(1)
$header  = $section.Headers.Item(1)
$header.Range.Find.Execute(
    "tsp",   $true, $false, $false, $false, $false,
    $true,   1,     $false, "tsi",  2) | Out-Null
